$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with new credentials/status
$ws.Range("A2").Value = "NinjaAlgo1"
$ws.Range("B2").Value = "@Algo3"
$ws.Range("C2").Value = "login"

# Add new row 3 with the previous credentials and a new status
$ws.Range("A3").Value = "NinjaAlgo"
$ws.Range("B3").Value = "@Algo123"
$ws.Range("C3").Value = "home"
